$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "43.219.87"
Set-TextValue "E2" "  -1.66%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.239.53"
Set-TextValue "E3" "  -1.17%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  -0.39%  "

# Row 5 - BNB
Set-TextValue "D5" "230.16"
Set-TextValue "E5" "  -0.28%  "

# Row 6 - XRP
Set-TextValue "D6" "0.643"
Set-TextValue "E6" "  +1.67%  "

# Row 7 - Solana
Set-TextValue "D7" "63.44"
Set-TextValue "E7" "  +0.73%  "

# Row 9 - Cardano
Set-TextValue "E9" "  +2.15%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.0948"
Set-TextValue "E10" "  -6.05%  "

# Row 11 - OKB
Set-TextValue "D11" "56.11"
Set-TextValue "E11" "  -0.36%  "

# Row 12 - Avalanche
Set-TextValue "D12" "27.31"
Set-TextValue "E12" "  +5.73%  "

# Row 13 - TRON
Set-TextValue "E13" "  -1.99%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue "D14" "2.568.03"
Set-TextValue "E14" "  -1.52%  "

# Row 15 - Chainlink
Set-TextValue "D15" "15.20"
Set-TextValue "E15" "  -2.75%  "

# Row 16 - Polkadot
Set-TextValue "D16" "6.02"
Set-TextValue "E16" "  +1.92%  "

# Row 17 - Polygon
Set-TextValue "D17" "0.822"
Set-TextValue "E17" "  +0.14%  "

# Row 18 - WrappedEther
Set-TextValue "D18" "2.238.73"
Set-TextValue "E18" "  -1.82%  "

# Row 19 - WrappedBTC
Set-TextValue "D19" "43.092.56"
Set-TextValue "E19" "  -1.91%  "

# Row 20 - ShibaInu
Set-TextValue "E20" "  -6.58%  "

# Row 21 - Litecoin
Set-TextValue "D21" "72.77"
Set-TextValue "E21" "  -1.07%  "

# Row 22 - Uniswap
Set-TextValue "D22" "6.05"
Set-TextValue "E22" "  +0.32%  "

# Row 23 - BitcoinCash
Set-TextValue "D23" "245.89"
Set-TextValue "E23" "  -4.05%  "

# Row 25 - WEMIXToken
Set-TextValue "D25" "3.72"
Set-TextValue "E25" "  +30.42%  "

# Row 26 - PancakeSwap
Set-TextValue "D26" "2.41"
Set-TextValue "E26" "  -1.94%  "

# Row 27 - Toncoin
Set-TextValue "D27" "2.22"
Set-TextValue "E27" "  -4.41%  "

# Row 28 - was Cosmos, now Monero
Set-TextValue "B28" "Monero"
Set-TextValue "C28" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D28" "173.43"
Set-TextValue "E28" "  +1.27%  "

# Row 29 - was Monero, now Cosmos
Set-TextValue "B29" "Cosmos"
Set-TextValue "C29" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D29" "9.67"
Set-TextValue "E29" "  -3.07%  "

# Row 30 - EthereumClassic
Set-TextValue "D30" "21.60"
Set-TextValue "E30" "  +3.52%  "

# Row 31 - was Kaspa, now ImmutableX
Set-TextValue "B31" "ImmutableX"
Set-TextValue "C31" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D31" "1.40"
Set-TextValue "E31" "  +0.69%  "

# Row 32 - was ImmutableX, now Kaspa
Set-TextValue "B32" "Kaspa"
Set-TextValue "C32" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D32" "0.128"
Set-TextValue "E32" "  -6.26%  "

# Row 33 - Stellar
Set-TextValue "E33" "  +0.58%  "

# Row 34 - Filecoin
Set-TextValue "E34" "  +4.64%  "

# Row 35 - Hedera
Set-TextValue "E35" "  -0.73%  "

# Row 36 - InternetComputer(DFINITY)
Set-TextValue "D36" "4.88"
Set-TextValue "E36" "  -1.65%  "

# Row 37 - RenderToken
Set-TextValue "D37" "3.57"
Set-TextValue "E37" "  -5.58%  "

# Row 38 - THORChain
Set-TextValue "D38" "6.28"
Set-TextValue "E38" "  -6.51%  "

# Row 39 - LidoDAOToken
Set-TextValue "E39" "  -2.92%  "

# Row 40 - VeChain
Set-TextValue "E40" "  -0.26%  "

# Row 41 - BinanceUSD
Set-TextValue "E41" "  -0.06%  "

# Row 42 - FraxShare
Set-TextValue "D42" "8.60"
Set-TextValue "E42" "  +1.91%  "

# Row 43 - FTXToken
Set-TextValue "D43" "4.42"
Set-TextValue "E43" "  +0.02%  "

# Row 44 - InjectiveProtocol
Set-TextValue "D44" "16.87"
Set-TextValue "E44" "  -2.62%  "

# Row 45 - Cronos
Set-TextValue "D45" "0.0940"
Set-TextValue "E45" "  -2.42%  "

# Row 46 - Aave
Set-TextValue "D46" "96.09"
Set-TextValue "E46" "  -1.26%  "

# Row 47 - TrustWalletToken
Set-TextValue "E47" "  -1.46%  "

# Row 48 - TerraClassic
Set-TextValue "D48" "0.000208"
Set-TextValue "E48" "  +0.23%  "

# Row 49 - Maker
Set-TextValue "D49" "1.430.67"
Set-TextValue "E49" "  -2.15%  "

# Row 50 - Celestia
Set-TextValue "E50" "  +4.30%  "

# Row 51 - HuobiToken
Set-TextValue "E51" "  +0.30%  "
